$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.4227110451726048
$ws.Range("C2").Value = 0.9311360074561618
$ws.Range("D2").Value = 0.5060024579641147
$ws.Range("G2").Value = 1.099153532933269
$ws.Range("H2").Value = 0.998

# Row 3
$ws.Range("B3").Value = 2.878893579773741
$ws.Range("G3").Value = 1.099153532933269
$ws.Range("H3").Value = 0.998

# Row 4
$ws.Range("B4").Value = 1.090838003737177
$ws.Range("C4").Value = 0.9461170227698321
$ws.Range("G4").Value = 1.099153532933269
$ws.Range("H4").Value = 0.998

# Row 5
$ws.Range("G5").Value = 1.099153532933269
$ws.Range("H5").Value = 0.998

# Row 6
$ws.Range("G6").Value = 1.099153532933269
$ws.Range("H6").Value = 0.998

# Row 7
$ws.Range("B7").Value = 1.271977756333029
$ws.Range("D7").Value = 0.8344729122263095
$ws.Range("G7").Value = 1.099153532933269
$ws.Range("H7").Value = 0.998

# Row 8
$ws.Range("G8").Value = 1.099153532933269
$ws.Range("H8").Value = 0.998

# Row 9
$ws.Range("B9").Value = 4.919241847624938
$ws.Range("C9").Value = 0.9411938629544762
$ws.Range("D9").Value = 1.690071404702129
$ws.Range("G9").Value = 1.099153532933269
$ws.Range("H9").Value = 0.998

# Row 10
$ws.Range("B10").Value = 0.8502471867908139
$ws.Range("C10").Value = 0.9974654394712144
$ws.Range("D10").Value = 0.7292237418674871
$ws.Range("G10").Value = 1.099153532933269
$ws.Range("H10").Value = 0.998

$wb.Save()
